$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "2024-05-11 21:26:21"
$ws.Range("B6").Value = 0.0002

$ws.Range("A7").Value = "2024-05-11 21:26:59"
$ws.Range("B7").Value = 0.001

$ws.Range("A8").Value = "2024-05-11 21:27:23"
$ws.Range("B8").Value = 0.0006000000000000001
